$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 87, shifting the existing rows 87-99 down to 89-101.
$ws.Rows("87:88").Insert()

# Row 87: new weekly entry (Especial)
$ws.Range("A87").Value = 1
$ws.Range("B87").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C87").Value = "Arica y Parinacota"
$ws.Range("D87").Value = 44491
$ws.Range("E87").Value = 15
$ws.Range("F87").Value = "Fruta"
$ws.Range("G87").Value = 100108
$ws.Range("H87").Value = "Tropicales y subtropicales"
$ws.Range("I87").Value = 100108002
$ws.Range("J87").Value = "Mango"
$ws.Range("K87").Value = "Sin especificar"
$ws.Range("L87").Value = "Especial"
$ws.Range("M87").Value = 456
$ws.Range("N87").Value = 4000
$ws.Range("O87").Value = 4500
$ws.Range("P87").Value = 4250
$ws.Range("Q87").Value = "$/bandeja 4 kilos"
$ws.Range("R87").Value = "Perú"
$ws.Range("S87").Value = 1062
$ws.Range("T87").Value = 4

# Row 88: new weekly entry (Primera)
$ws.Range("A88").Value = 1
$ws.Range("B88").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C88").Value = "Arica y Parinacota"
$ws.Range("D88").Value = 44491
$ws.Range("E88").Value = 15
$ws.Range("F88").Value = "Fruta"
$ws.Range("G88").Value = 100108
$ws.Range("H88").Value = "Tropicales y subtropicales"
$ws.Range("I88").Value = 100108002
$ws.Range("J88").Value = "Mango"
$ws.Range("K88").Value = "Sin especificar"
$ws.Range("L88").Value = "Primera"
$ws.Range("M88").Value = 456
$ws.Range("N88").Value = 4000
$ws.Range("O88").Value = 4500
$ws.Range("P88").Value = 4250
$ws.Range("Q88").Value = "$/bandeja 4 kilos"
$ws.Range("R88").Value = "Perú"
$ws.Range("S88").Value = 1062
$ws.Range("T88").Value = 4
